# updated r1 to w1
# The "race" row (row 4) encodes its categorical codes as "r1, b2, as3, ain4, r0"
# (allowed_range, column F) with a matching legend "r1=White, b2=Black, ..."
# (notes, column G). The code for White is being renamed from r1 to w1 in both
# places.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G4").Value = "w1=White, b2=Black, as3=Asian/PI, ain4=Native American, r0=Unknown"
$ws.Range("F4").Value = "w1, b2, as3, ain4, r0"

# Move the active selection, matching the post-edit cursor position.
$ws.Range("E23").Select()
